# Append the new form-submission row (row 8) to the sheet's data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Nara"
$ws.Range("B8").Value = "nareay@cat.com"
$ws.Range("C8").Value = "Employee"
$ws.Range("D8").Value = "2025-10-01T18:37:43.378Z"
